# Weekly price-list update: insert a new record at the top of the data
# block (row 78) for "Vega Modelo de Temuco - Espinaca", pushing the
# existing rows 78-95 down to 79-96.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 78; this also shifts
# rows 78-95 down to 79-96 and extends the used range to row 96.
$ws.Rows.Item(78).Insert()

# Populate the new row 78 with this week's record. The "constant"
# columns (A,B,C,E,F,G,H,I,N,O,Q,R) repeat the same values used by
# every other row of this data set.
$ws.Range("A78").Value = 10
$ws.Range("B78").Value = "Vega Modelo de Temuco"
$ws.Range("C78").Value = "La Araucanía"
$ws.Range("D78").Value = 44508
$ws.Range("E78").Value = 9
$ws.Range("F78").Value = 100112012
$ws.Range("G78").Value = "Espinaca"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 50
$ws.Range("K78").Value = 8000
$ws.Range("L78").Value = 8000
$ws.Range("M78").Value = 8000
$ws.Range("N78").Value = "$/docena de atados"
$ws.Range("O78").Value = "Región de La Araucanía"
$ws.Range("P78").Value = 2667
$ws.Range("Q78").Value = 3
$ws.Range("R78").Value = "Hortaliza"
